$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.23539999999999
$ws.Range("C12").Value = -11.3212
$ws.Range("E14").Value = 16.9894
$ws.Range("E26").Value = 16.37359999999999
$ws.Range("E31").Value = 16.8349
$ws.Range("C32").Value = -13.60450000000001
$ws.Range("E35").Value = 16.46710000000001
$ws.Range("C36").Value = -13.06680000000002
$ws.Range("E37").Value = 16.69290000000001
$ws.Range("C38").Value = -13.08930000000001
$ws.Range("E45").Value = 16.5929
$ws.Range("C46").Value = -14.22819999999999
$ws.Range("C54").Value = -12.7228
$ws.Range("C55").Value = -13.43839999999999
$ws.Range("E57").Value = 16.8133
$ws.Range("C67").Value = -10.81470000000001
$ws.Range("C69").Value = -12.42429999999999
$ws.Range("C72").Value = -11.5643
$ws.Range("C91").Value = -10.5483
$ws.Range("C99").Value = -13.4238
$ws.Range("E100").Value = 16.4876
$ws.Range("E102").Value = 16.70829999999999
